# .hide() aitab asju ära peita
# Add the Valik1 / Valik2 / Valik3 / Vanusegrupp columns to the Grupp_3 sheet
# (sheet3), matching the layout already used on Grupp_1 / Grupp_2, and record
# the age group (3) for every respondent row on that sheet.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item(3)

# New header cells (C1:E1) - reuse the exact text already used as headers on
# the other two sheets so the shared-string table is not fattened, and move
# the "Vastus" header out to its new home in column F.
$ws3.Range("C1").Value = "Valik1"
$ws3.Range("D1").Value = "Valik2"
$ws3.Range("E1").Value = "Valik3"
$ws3.Range("F1").Value = "Vanusegrupp"

# Every respondent in Grupp_3 belongs to age group 3.
$ws3.Range("F2:F7").Value = 3

# Leave the sheet's selection where the author left it.
$ws3.Range("C2").Select() | Out-Null

# The workbook was last saved with Grupp_2 as the active tab.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate() | Out-Null
$ws2.Range("C11").Select() | Out-Null
